# Auto-generated Excel COM-interop script
# Applies numeric value updates to ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 2752.4075
$ws.Cells.Item(15, 9).Value = 2752.4075
$ws.Cells.Item(15, 11).Value = 8257.2225
$ws.Cells.Item(15, 13).Value = -8088.2225

$ws.Cells.Item(18, 8).Value = 10968.308
$ws.Cells.Item(18, 9).Value = 5689.909
$ws.Cells.Item(18, 10).Value = 39999.5
$ws.Cells.Item(18, 11).Value = 5689.909
$ws.Cells.Item(18, 12).Value = 39999.5
$ws.Cells.Item(18, 13).Value = -5405.909
$ws.Cells.Item(18, 14).Value = -40567.5

$ws.Cells.Item(86, 8).Value = 5792.3447
$ws.Cells.Item(86, 9).Value = 5780.421
$ws.Cells.Item(86, 10).Value = 5815
$ws.Cells.Item(86, 11).Value = 5780.421
$ws.Cells.Item(86, 12).Value = 5815
$ws.Cells.Item(86, 13).Value = -4657.421
$ws.Cells.Item(86, 14).Value = -8061

$ws.Cells.Item(89, 8).Value = 5792.3447
$ws.Cells.Item(89, 9).Value = 5780.421
$ws.Cells.Item(89, 10).Value = 5815
$ws.Cells.Item(89, 11).Value = 28902.105
$ws.Cells.Item(89, 12).Value = 29075
$ws.Cells.Item(89, 13).Value = -23286.105
$ws.Cells.Item(89, 14).Value = -40307

$ws.Cells.Item(92, 8).Value = 866.3871
$ws.Cells.Item(92, 9).Value = 479.21054
$ws.Cells.Item(92, 10).Value = 1479.4166
$ws.Cells.Item(92, 11).Value = 479.21054
$ws.Cells.Item(92, 12).Value = 1479.4166
$ws.Cells.Item(92, 13).Value = 768.78946
$ws.Cells.Item(92, 14).Value = -3975.4166

$ws.Cells.Item(132, 8).Value = 1913.4193
$ws.Cells.Item(132, 9).Value = 1838.5862
$ws.Cells.Item(132, 11).Value = 5515.7586
$ws.Cells.Item(132, 13).Value = -2985.7586

$ws.Cells.Item(135, 8).Value = 7712.778
$ws.Cells.Item(135, 9).Value = 2217.6924
$ws.Cells.Item(135, 10).Value = 22000
$ws.Cells.Item(135, 11).Value = 19959.2316
$ws.Cells.Item(135, 12).Value = 198000
$ws.Cells.Item(135, 13).Value = -17424.2316
$ws.Cells.Item(135, 14).Value = -203070

$ws.Cells.Item(137, 8).Value = 9565.31
$ws.Cells.Item(137, 9).Value = 4596.7393
$ws.Cells.Item(137, 10).Value = 13136.469
$ws.Cells.Item(137, 11).Value = 13790.2179
$ws.Cells.Item(137, 12).Value = 39409.407
$ws.Cells.Item(137, 13).Value = -11240.2179
$ws.Cells.Item(137, 14).Value = -44509.407

$ws.Cells.Item(138, 8).Value = 5803.7646
$ws.Cells.Item(138, 10).Value = 5677.195
$ws.Cells.Item(138, 12).Value = 17031.585
$ws.Cells.Item(138, 14).Value = -27311.585

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11127075
$ws.Cells.Item(32, 9).Value = 13524919
$ws.Cells.Item(32, 10).Value = 37046
$ws.Cells.Item(32, 11).Value = 13524919
$ws.Cells.Item(32, 12).Value = 37046
$ws.Cells.Item(32, 13).Value = -13524632
$ws.Cells.Item(32, 14).Value = -37620

$ws.Cells.Item(45, 8).Value = 2328.2856
$ws.Cells.Item(45, 9).Value = 2127.875
$ws.Cells.Item(45, 10).Value = 2451.6155
$ws.Cells.Item(45, 11).Value = 2127.875
$ws.Cells.Item(45, 12).Value = 2451.6155
$ws.Cells.Item(45, 13).Value = -1750.875
$ws.Cells.Item(45, 14).Value = -3205.6155

$ws.Cells.Item(61, 8).Value = 16703802
$ws.Cells.Item(61, 9).Value = 22729764
$ws.Cells.Item(61, 10).Value = 132406.88
$ws.Cells.Item(61, 11).Value = 22729764
$ws.Cells.Item(61, 12).Value = 132406.88
$ws.Cells.Item(61, 13).Value = -22729552
$ws.Cells.Item(61, 14).Value = -132830.88

$ws.Cells.Item(74, 8).Value = 5958205
$ws.Cells.Item(74, 9).Value = 8930170
$ws.Cells.Item(74, 10).Value = 14273.786
$ws.Cells.Item(74, 11).Value = 8930170
$ws.Cells.Item(74, 12).Value = 14273.786
$ws.Cells.Item(74, 13).Value = -8929296
$ws.Cells.Item(74, 14).Value = -16021.786

$ws.Cells.Item(77, 8).Value = 5958205
$ws.Cells.Item(77, 9).Value = 8930170
$ws.Cells.Item(77, 10).Value = 14273.786
$ws.Cells.Item(77, 11).Value = 44650850
$ws.Cells.Item(77, 12).Value = 71368.93
$ws.Cells.Item(77, 13).Value = -44646482
$ws.Cells.Item(77, 14).Value = -80104.93

$ws.Cells.Item(136, 8).Value = 16703802
$ws.Cells.Item(136, 9).Value = 22729764
$ws.Cells.Item(136, 10).Value = 132406.88
$ws.Cells.Item(136, 11).Value = 68189292
$ws.Cells.Item(136, 12).Value = 397220.64
$ws.Cells.Item(136, 13).Value = -68186742
$ws.Cells.Item(136, 14).Value = -402320.64

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 27360.834
$ws.Cells.Item(134, 9).Value = 3369.2
$ws.Cells.Item(134, 10).Value = 87339.914
$ws.Cells.Item(134, 11).Value = 10107.6
$ws.Cells.Item(134, 12).Value = 262019.742
$ws.Cells.Item(134, 13).Value = -7572.599999999999
$ws.Cells.Item(134, 14).Value = -267089.742

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 1626.5
$ws.Cells.Item(7, 10).Value = 3487.077
$ws.Cells.Item(7, 12).Value = 3487.077
$ws.Cells.Item(7, 14).Value = -3713.077

$ws.Cells.Item(31, 8).Value = 369651.22
$ws.Cells.Item(31, 9).Value = 10939.833
$ws.Cells.Item(31, 10).Value = 510016.53
$ws.Cells.Item(31, 11).Value = 10939.833
$ws.Cells.Item(31, 12).Value = 510016.53
$ws.Cells.Item(31, 13).Value = -10644.833
$ws.Cells.Item(31, 14).Value = -510606.53

$ws.Cells.Item(34, 8).Value = 369651.22
$ws.Cells.Item(34, 9).Value = 10939.833
$ws.Cells.Item(34, 10).Value = 510016.53
$ws.Cells.Item(34, 11).Value = 10939.833
$ws.Cells.Item(34, 12).Value = 510016.53
$ws.Cells.Item(34, 13).Value = -10737.833
$ws.Cells.Item(34, 14).Value = -510420.53

$ws.Cells.Item(41, 8).Value = 2000
$ws.Cells.Item(41, 9).Value = 2000
$ws.Cells.Item(41, 11).Value = 2000
$ws.Cells.Item(41, 13).Value = -1572

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 3395.0466
$ws.Cells.Item(68, 9).Value = 3056.9167
$ws.Cells.Item(68, 10).Value = 3525.9355
$ws.Cells.Item(68, 11).Value = 9170.750100000001
$ws.Cells.Item(68, 12).Value = 10577.8065
$ws.Cells.Item(68, 13).Value = -8359.750100000001
$ws.Cells.Item(68, 14).Value = -12199.8065

$ws.Cells.Item(71, 8).Value = 3395.0466
$ws.Cells.Item(71, 9).Value = 3056.9167
$ws.Cells.Item(71, 10).Value = 3525.9355
$ws.Cells.Item(71, 11).Value = 27512.2503
$ws.Cells.Item(71, 12).Value = 31733.4195
$ws.Cells.Item(71, 13).Value = -23456.2503
$ws.Cells.Item(71, 14).Value = -39845.4195

$ws.Cells.Item(107, 8).Value = 1330.2222
$ws.Cells.Item(107, 9).Value = 999.3333
$ws.Cells.Item(107, 11).Value = 2997.9999
$ws.Cells.Item(107, 13).Value = -1077.9999

$ws.Cells.Item(113, 8).Value = 1444.6
$ws.Cells.Item(113, 10).Value = 1352.0454
$ws.Cells.Item(113, 12).Value = 4056.1362
$ws.Cells.Item(113, 14).Value = -8396.1362

$ws.Cells.Item(132, 8).Value = 2244.7144
$ws.Cells.Item(132, 9).Value = 2216.2666
$ws.Cells.Item(132, 10).Value = 2315.8333
$ws.Cells.Item(132, 11).Value = 19946.3994
$ws.Cells.Item(132, 12).Value = 20842.4997
$ws.Cells.Item(132, 13).Value = -17416.3994
$ws.Cells.Item(132, 14).Value = -25902.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8667.833
$ws.Cells.Item(70, 9).Value = 8001.75
$ws.Cells.Item(70, 11).Value = 8001.75
$ws.Cells.Item(70, 13).Value = -7731.75

$ws.Cells.Item(73, 8).Value = 8667.833
$ws.Cells.Item(73, 9).Value = 8001.75
$ws.Cells.Item(73, 11).Value = 8001.75
$ws.Cells.Item(73, 13).Value = -7065.75

$ws.Cells.Item(113, 8).Value = 3885.2727
$ws.Cells.Item(113, 9).Value = 2945.4
$ws.Cells.Item(113, 10).Value = 4668.5
$ws.Cells.Item(113, 11).Value = 2945.4
$ws.Cells.Item(113, 12).Value = 4668.5
$ws.Cells.Item(113, 13).Value = -775.4000000000001
$ws.Cells.Item(113, 14).Value = -9008.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 40026.94
$ws.Cells.Item(136, 9).Value = 7190.4375
$ws.Cells.Item(136, 11).Value = 21571.3125
$ws.Cells.Item(136, 13).Value = -19021.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2865.7
$ws.Cells.Item(132, 9).Value = 2220.484
$ws.Cells.Item(132, 10).Value = 5088.1113
$ws.Cells.Item(132, 11).Value = 6661.451999999999
$ws.Cells.Item(132, 12).Value = 15264.3339
$ws.Cells.Item(132, 13).Value = -4131.451999999999
$ws.Cells.Item(132, 14).Value = -20324.3339

$ws.Cells.Item(136, 8).Value = 2200.9167
$ws.Cells.Item(136, 9).Value = 1540.7
$ws.Cells.Item(136, 11).Value = 4622.1
$ws.Cells.Item(136, 13).Value = -2072.1
